$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.901.84'
$ws.Range('E2').Value = '  +2.66%  '
$ws.Range('D3').Value = '1.874.15'
$ws.Range('E4').Value = '  -0.58%  '
$ws.Range('D5').Value = "'314.15"
$ws.Range('E5').Value = '  +1.17%  '
$ws.Range('D6').Value = "'1.013"
$ws.Range('E6').Value = '  -0.32%  '
$ws.Range('D7').Value = "'0.4839"
$ws.Range('E7').Value = '  +1.30%  '
$ws.Range('D8').Value = "'0.3828"
$ws.Range('E8').Value = '  +3.46%  '
$ws.Range('E9').Value = '  +1.51%  '
$ws.Range('D10').Value = "'0.9427"
$ws.Range('E10').Value = '  +0.90%  '
$ws.Range('D11').Value = "'21.11"
$ws.Range('E11').Value = '  +5.82%  '
$ws.Range('D12').Value = "'0.07820"
$ws.Range('E12').Value = '  +0.02%  '
$ws.Range('D13').Value = '1.882.45'
$ws.Range('E13').Value = '  +1.34%  '
$ws.Range('D14').Value = "'5.499"
$ws.Range('E14').Value = '  +1.86%  '
$ws.Range('D15').Value = "'6.621"
$ws.Range('E15').Value = '  +1.68%  '
$ws.Range('D16').Value = "'91.33"
$ws.Range('E16').Value = '  +1.75%  '
$ws.Range('E17').Value = '  -0.58%  '
$ws.Range('D18').Value = "'0.000008894"
$ws.Range('E18').Value = '  +2.05%  '
$ws.Range('E19').Value = '  -0.54%  '
$ws.Range('D20').Value = '27.917.66'
$ws.Range('E20').Value = '  +2.62%  '
$ws.Range('D21').Value = "'14.89"
$ws.Range('E21').Value = '  +1.67%  '
$ws.Range('D22').Value = "'5.132"
$ws.Range('E22').Value = '  +0.94%  '
$ws.Range('D23').Value = '2.114.79'
$ws.Range('E23').Value = '  +1.43%  '
$ws.Range('E24').Value = '  +1.72%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').Value = "'1.948"
$ws.Range('E25').Value = '  +0.10%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = "'157.18"
$ws.Range('E26').Value = '  +2.52%  '
$ws.Range('D27').Value = "'18.60"
$ws.Range('E27').Value = '  +0.95%  '
$ws.Range('D28').Value = "'2.063"
$ws.Range('E28').Value = '  +3.28%  '
$ws.Range('D29').Value = "'116.14"
$ws.Range('E29').Value = '  +0.89%  '
$ws.Range('E30').Value = '  +1.28%  '
$ws.Range('D31').Value = "'0.08932"
$ws.Range('E31').Value = '  +0.64%  '
$ws.Range('D32').Value = "'3.329"
$ws.Range('E32').Value = '  +0.58%  '
$ws.Range('D33').Value = "'1.233"
$ws.Range('E33').Value = '  +4.27%  '
$ws.Range('D34').Value = "'0.7711"
$ws.Range('E34').Value = '  +4.49%  '
$ws.Range('E35').Value = '  +2.55%  '
$ws.Range('D36').Value = "'2.710"
$ws.Range('E36').Value = '  +0.69%  '
$ws.Range('E37').Value = '  +1.50%  '
$ws.Range('D38').Value = "'0.02052"
$ws.Range('E38').Value = '  +2.73%  '
$ws.Range('D39').Value = "'0.5634"
$ws.Range('E39').Value = '  +6.25%  '
$ws.Range('D40').Value = "'0.05372"
$ws.Range('E40').Value = '  +2.20%  '
$ws.Range('D41').Value = "'2.997"
$ws.Range('E41').Value = '  +0.64%  '
$ws.Range('D42').Value = "'7.049"
$ws.Range('E42').Value = '  -0.04%  '
$ws.Range('E43').Value = '  +3.21%  '
$ws.Range('E44').Value = '  +0.70%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = "'10.76"
$ws.Range('E45').Value = '  +1.24%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').Value = "'0.4898"
$ws.Range('E46').Value = '  +3.03%  '
$ws.Range('D47').Value = "'105.47"
$ws.Range('E47').Value = '  +3.09%  '
$ws.Range('E48').Value = '  -0.56%  '
$ws.Range('D49').Value = "'1.671"
$ws.Range('E49').Value = '  +2.67%  '
$ws.Range('D50').Value = "'68.13"
$ws.Range('E50').Value = '  +3.05%  '
$ws.Range('D51').Value = "'0.06118"
$ws.Range('E51').Value = '  +0.86%  '
